$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value2 = 709.8
$ws.Range("I103").Value2 = 1149.5
$ws.Range("J103").Value2 = 416.66666
$ws.Range("K103").Value2 = 3448.5
$ws.Range("L103").Value2 = 1249.99998
$ws.Range("M103").Value2 = -2862.5
$ws.Range("N103").Value2 = -2421.99998
$ws.Range("H116").Value2 = 8246.565000000001
$ws.Range("I116").Value2 = 3633.4119
$ws.Range("K116").Value2 = 3633.4119
$ws.Range("M116").Value2 = -191.4119000000001
$ws.Range("H118").Value2 = 746
$ws.Range("I118").Value2 = 592.1111
$ws.Range("J118").Value2 = 1207.6666
$ws.Range("K118").Value2 = 1776.3333
$ws.Range("L118").Value2 = 3622.9998
$ws.Range("M118").Value2 = -119.3332999999998
$ws.Range("N118").Value2 = -6936.9998
$ws.Range("H135").Value2 = 690860.3
$ws.Range("J135").Value2 = 1317.2
$ws.Range("L135").Value2 = 11854.8
$ws.Range("N135").Value2 = -16924.8
$ws.Range("H141").Value2 = 2978.4285
$ws.Range("J141").Value2 = 3063.3333
$ws.Range("L141").Value2 = 9189.999899999999
$ws.Range("N141").Value2 = -19549.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5250.156
$ws.Range("I32").Value2 = 5250.156
$ws.Range("K32").Value2 = 5250.156
$ws.Range("M32").Value2 = -4963.156
$ws.Range("H61").Value2 = 2865.8076
$ws.Range("I61").Value2 = 2348.348
$ws.Range("K61").Value2 = 2348.348
$ws.Range("M61").Value2 = -2136.348
$ws.Range("H74").Value2 = 2253.5
$ws.Range("I74").Value2 = 1861.2142
$ws.Range("J74").Value2 = 4999.5
$ws.Range("K74").Value2 = 1861.2142
$ws.Range("L74").Value2 = 4999.5
$ws.Range("M74").Value2 = -987.2141999999999
$ws.Range("N74").Value2 = -6747.5
$ws.Range("H77").Value2 = 2253.5
$ws.Range("I77").Value2 = 1861.2142
$ws.Range("J77").Value2 = 4999.5
$ws.Range("K77").Value2 = 9306.071
$ws.Range("L77").Value2 = 24997.5
$ws.Range("M77").Value2 = -4938.071
$ws.Range("N77").Value2 = -33733.5
$ws.Range("H110").Value2 = 126275.15
$ws.Range("I110").Value2 = 132865.62
$ws.Range("J110").Value2 = 1056
$ws.Range("K110").Value2 = 132865.62
$ws.Range("L110").Value2 = 1056
$ws.Range("M110").Value2 = -130820.62
$ws.Range("N110").Value2 = -5146
$ws.Range("H132").Value2 = 2462.5405
$ws.Range("I132").Value2 = 1818.2258
$ws.Range("J132").Value2 = 5791.5
$ws.Range("K132").Value2 = 5454.6774
$ws.Range("L132").Value2 = 17374.5
$ws.Range("M132").Value2 = -2924.6774
$ws.Range("N132").Value2 = -22434.5
$ws.Range("H136").Value2 = 2865.8076
$ws.Range("I136").Value2 = 2348.348
$ws.Range("K136").Value2 = 7045.044
$ws.Range("M136").Value2 = -4495.044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value2 = 3704
$ws.Range("J24").Value2 = 2000
$ws.Range("L24").Value2 = 2000
$ws.Range("N24").Value2 = -2470
$ws.Range("H74").Value2 = 97317.39999999999
$ws.Range("J74").Value2 = 97317.39999999999
$ws.Range("L74").Value2 = 97317.39999999999
$ws.Range("N74").Value2 = -99189.39999999999
$ws.Range("H77").Value2 = 97317.39999999999
$ws.Range("J77").Value2 = 97317.39999999999
$ws.Range("L77").Value2 = 291952.2
$ws.Range("N77").Value2 = -301312.2
$ws.Range("H81").Value2 = 82996
$ws.Range("J81").Value2 = 82996
$ws.Range("L81").Value2 = 82996
$ws.Range("N81").Value2 = -85118
$ws.Range("H84").Value2 = 82996
$ws.Range("J84").Value2 = 82996
$ws.Range("L84").Value2 = 248988
$ws.Range("N84").Value2 = -259596
$ws.Range("H99").Value2 = 1071.5714
$ws.Range("I99").Value2 = 1088.25
$ws.Range("K99").Value2 = 1088.25
$ws.Range("M99").Value2 = 409.75
$ws.Range("H134").Value2 = 29262.578
$ws.Range("I134").Value2 = 1935.9333
$ws.Range("K134").Value2 = 5807.7999
$ws.Range("M134").Value2 = -3272.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1036.0769
$ws.Range("I16").Value2 = 832.6667
$ws.Range("K16").Value2 = 832.6667
$ws.Range("M16").Value2 = -545.6667
$ws.Range("H58").Value2 = 282599.38
$ws.Range("I58").Value2 = 404130.56
$ws.Range("J58").Value2 = 6392.091
$ws.Range("K58").Value2 = 404130.56
$ws.Range("L58").Value2 = 6392.091
$ws.Range("M58").Value2 = -403927.56
$ws.Range("N58").Value2 = -6798.091
$ws.Range("H113").Value2 = 1036.0769
$ws.Range("I113").Value2 = 832.6667
$ws.Range("K113").Value2 = 832.6667
$ws.Range("M113").Value2 = 1337.3333
$ws.Range("H134").Value2 = 670002.1
$ws.Range("I134").Value2 = 402881.97
$ws.Range("K134").Value2 = 1208645.91
$ws.Range("M134").Value2 = -1206110.91
$ws.Range("H136").Value2 = 282599.38
$ws.Range("I136").Value2 = 404130.56
$ws.Range("J136").Value2 = 6392.091
$ws.Range("K136").Value2 = 1212391.68
$ws.Range("L136").Value2 = 19176.273
$ws.Range("M136").Value2 = -1209841.68
$ws.Range("N136").Value2 = -24276.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 1642356.9
$ws.Range("I4").Value2 = 2067101.5
$ws.Range("K4").Value2 = 6201304.5
$ws.Range("M4").Value2 = -6201192.5
$ws.Range("H5").Value2 = 34249.625
$ws.Range("I5").Value2 = 53634.4
$ws.Range("J5").Value2 = 1941.6666
$ws.Range("K5").Value2 = 160903.2
$ws.Range("L5").Value2 = 5824.9998
$ws.Range("M5").Value2 = -160791.2
$ws.Range("N5").Value2 = -6048.9998
$ws.Range("H11").Value2 = 625497.25
$ws.Range("I11").Value2 = 1111348.4
$ws.Range("J11").Value2 = 831.5
$ws.Range("K11").Value2 = 3334045.2
$ws.Range("L11").Value2 = 2494.5
$ws.Range("M11").Value2 = -3333905.2
$ws.Range("N11").Value2 = -2774.5
$ws.Range("H107").Value2 = 38263.5
$ws.Range("I107").Value2 = 879.75
$ws.Range("J107").Value2 = 66301.31
$ws.Range("K107").Value2 = 2639.25
$ws.Range("L107").Value2 = 198903.93
$ws.Range("M107").Value2 = -719.25
$ws.Range("N107").Value2 = -202743.93
$ws.Range("H109").Value2 = 92390.27
$ws.Range("I109").Value2 = 715.5
$ws.Range("K109").Value2 = 2146.5
$ws.Range("M109").Value2 = -1106.5
$ws.Range("H120").Value2 = 6617.5
$ws.Range("J120").Value2 = 2000
$ws.Range("L120").Value2 = 6000
$ws.Range("N120").Value2 = -15676
$ws.Range("H135").Value2 = 34249.625
$ws.Range("I135").Value2 = 53634.4
$ws.Range("J135").Value2 = 1941.6666
$ws.Range("K135").Value2 = 482709.6
$ws.Range("L135").Value2 = 17474.9994
$ws.Range("M135").Value2 = -480174.6
$ws.Range("N135").Value2 = -22544.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 38467890
$ws.Range("I70").Value2 = 7210.278
$ws.Range("J70").Value2 = 125004420
$ws.Range("K70").Value2 = 7210.278
$ws.Range("L70").Value2 = 125004420
$ws.Range("M70").Value2 = -6940.278
$ws.Range("N70").Value2 = -125004960
$ws.Range("H73").Value2 = 38467890
$ws.Range("I73").Value2 = 7210.278
$ws.Range("J73").Value2 = 125004420
$ws.Range("K73").Value2 = 7210.278
$ws.Range("L73").Value2 = 125004420
$ws.Range("M73").Value2 = -6274.278
$ws.Range("N73").Value2 = -125006292
$ws.Range("H113").Value2 = 9321.944
$ws.Range("I113").Value2 = 3407
$ws.Range("K113").Value2 = 3407
$ws.Range("M113").Value2 = -1237
$ws.Range("H132").Value2 = 397457.22
$ws.Range("I132").Value2 = 592252.3
$ws.Range("K132").Value2 = 1776756.9
$ws.Range("M132").Value2 = -1774226.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 375407.62
$ws.Range("I7").Value2 = 5256.357
$ws.Range("K7").Value2 = 5256.357
$ws.Range("M7").Value2 = -5144.357
$ws.Range("H40").Value2 = 2503347.5
$ws.Range("I40").Value2 = 3336171.8
$ws.Range("J40").Value2 = 4874.8
$ws.Range("K40").Value2 = 3336171.8
$ws.Range("L40").Value2 = 4874.8
$ws.Range("M40").Value2 = -3336035.8
$ws.Range("N40").Value2 = -5146.8
$ws.Range("H126").Value2 = 375407.62
$ws.Range("I126").Value2 = 5256.357
$ws.Range("K126").Value2 = 15769.071
$ws.Range("M126").Value2 = -13299.071
$ws.Range("H132").Value2 = 5056.25
$ws.Range("H136").Value2 = 1545836
$ws.Range("I136").Value2 = 2505870
$ws.Range("K136").Value2 = 7517610
$ws.Range("M136").Value2 = -7515060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value2 = 1000
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 1000
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 1000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value2 = -1226
$ws.Range("I74").Value2 = 4743.5
$ws.Range("J74").Value2 = 7287.8887
$ws.Range("K74").Value2 = 4743.5
$ws.Range("L74").Value2 = 7287.8887
$ws.Range("M74").Value2 = -3807.5
$ws.Range("N74").Value2 = -9159.8887
$ws.Range("I77").Value2 = 4743.5
$ws.Range("J77").Value2 = 7287.8887
$ws.Range("K77").Value2 = 14230.5
$ws.Range("L77").Value2 = 21863.6661
$ws.Range("M77").Value2 = -9550.5
$ws.Range("N77").Value2 = -31223.6661
$ws.Range("H122").Value2 = 24392506
$ws.Range("I122").Value2 = 26317638
$ws.Range("K122").Value2 = 78952914
$ws.Range("M122").Value2 = -78950464
$ws.Range("H126").Value2 = 2505
$ws.Range("I126").Value2 = 581.5714
$ws.Range("K126").Value2 = 1744.7142
$ws.Range("M126").Value2 = 725.2857999999999
$ws.Range("H132").Value2 = 19327.896
$ws.Range("I132").Value2 = 1526.5098
$ws.Range("K132").Value2 = 4579.5294
$ws.Range("M132").Value2 = -2049.5294
